# Apply 2022-05-16 data update to "Fonds de solidarite" VOLET1 regional/classe effectif dataset.
# For each affected row, update nombre_aides (C), nombre_entreprises (D, only row 184) and montant_total (E).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 64;  C = 5203;   D = $null;  E = 20360065 },
    @{ Row = 66;  C = 768;    D = $null;  E = 9929275 },
    @{ Row = 70;  C = 15726;  D = $null;  E = 24658731 },
    @{ Row = 73;  C = 2456;   D = $null;  E = 7383072 },
    @{ Row = 74;  C = 948;    D = $null;  E = 4258463 },
    @{ Row = 75;  C = 403;    D = $null;  E = 2849572 },
    @{ Row = 76;  C = 128;    D = $null;  E = 2483190 },
    @{ Row = 91;  C = 151098; D = $null;  E = 482084331 },
    @{ Row = 92;  C = 409006; D = $null;  E = 1593616902 },
    @{ Row = 93;  C = 209490; D = $null;  E = 1307691204 },
    @{ Row = 94;  C = 94144;  D = $null;  E = 915770565 },
    @{ Row = 95;  C = 50722;  D = $null;  E = 930422875 },
    @{ Row = 98;  C = 810;    D = $null;  E = 117791167 },
    @{ Row = 101; C = 179;    D = $null;  E = 32052522 },
    @{ Row = 104; C = 135216; D = $null;  E = 272104589 },
    @{ Row = 105; C = 8170;   D = $null;  E = 16872184 },
    @{ Row = 107; C = 6390;   D = $null;  E = 21955635 },
    @{ Row = 114; C = 3798;   D = $null;  E = 9098224 },
    @{ Row = 116; C = 4558;   D = $null;  E = 20484785 },
    @{ Row = 117; C = 1913;   D = $null;  E = 12352831 },
    @{ Row = 118; C = 976;    D = $null;  E = 11793478 },
    @{ Row = 122; C = 8488;   D = $null;  E = 12672414 },
    @{ Row = 132; C = 30288;  D = $null;  E = 174222945 },
    @{ Row = 144; C = 24413;  D = $null;  E = 201797181 },
    @{ Row = 173; C = 96858;  D = $null;  E = 327924730 },
    @{ Row = 184; C = 68734;  D = 13881;  E = 134170601 }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Cells.Item($r, 3).Value = $u.C
    if ($null -ne $u.D) {
        $ws.Cells.Item($r, 4).Value = $u.D
    }
    $ws.Cells.Item($r, 5).Value = $u.E
}
